$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.810.75'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  +7.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.504.39'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +7.86%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.31'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +13.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '556.79'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +8.31%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.614'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  +3.75%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.493.41'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +7.84%  '
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.641'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +7.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.68'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +4.84%  '
$ws.Range('E12').Value = '  +15.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000276'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +9.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.50'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +6.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.068.28'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +7.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.503.92'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +7.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.160.23'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +7.60%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.121'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +5.51%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.38'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +7.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.90'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +9.81%  '
$ws.Range('E21').Value = '  +7.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '405.43'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +10.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.07'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +12.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.96'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +7.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.99'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +7.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.21'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +9.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.93'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +11.75%  '
$ws.Range('E28').Value = '  +2.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.90'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +7.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.66'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +6.70%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '695.71'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +10.15%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '30.45'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +7.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.91'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +6.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.74'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +6.29%  '
$ws.Range('E35').Value = '  +8.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '60.96'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +4.44%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.07'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +8.76%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0829'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +22.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.405'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +7.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.40'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +26.07%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.135'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +12.20%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.81'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +17.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.064.61'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +6.40%  '
$ws.Range('E46').Value = '  +13.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0422'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +9.13%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.79'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +6.11%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.22'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +8.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.85'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +15.80%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.131'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +6.41%  '
